# Prepend "Design: " to the feedback bullet answers in the table.
# Each Find/Execute targets the exact original text of a run's first
# sentence/paragraph and replaces it with the "Design: " prefixed version.
# MatchCase=$true, MatchWholeWord=$false, Wrap=wdFindContinue(1),
# Format=$false, Replace=wdReplaceAll(2).

$d = $word.ActiveDocument

function Add-DesignPrefix([string]$oldText, [string]$newText) {
    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $newText, 2) | Out-Null
}

Add-DesignPrefix "Mostly" "Design: Mostly"

Add-DesignPrefix "Kymys ei täysin ymmärretävä, tarkoitetaanko tällä sitä miten sähkö on onnistunut?" `
                  "Design: Kymys ei täysin ymmärretävä, tarkoitetaanko tällä sitä miten sähkö on onnistunut?"

Add-DesignPrefix "Nothing special to mention" "Design: Nothing special to mention"

Add-DesignPrefix "Omien järjestelmien osalta (sähkönjakelu) vaikuttaa moneen järjestelmään säännöt esim. tarpeet tuplasyötöistä, syötönvaihdot jne." `
                  "Design: Omien järjestelmien osalta (sähkönjakelu) vaikuttaa moneen järjestelmään säännöt esim. tarpeet tuplasyötöistä, syötönvaihdot jne."

Add-DesignPrefix "91xx Went well, some minor budget challenges" "Design: 91xx Went well, some minor budget challenges"

Add-DesignPrefix "tarjouspyynnöt ja tarjouskierros monimutkainen ja aikaa vievä prosessi. Kaikkiaan kun saadaan 3 tarjousta sen jälkeen alkaa uusi tinkauskierros jossa taas hintaa lasketaan. Kun on päästy loppusuoralle asiat on jo muuttuneet. (Esim. tullut uusia sähkönkuluttajia tai uutta tietoa). " `
                  "Design: tarjouspyynnöt ja tarjouskierros monimutkainen ja aikaa vievä prosessi. Kaikkiaan kun saadaan 3 tarjousta sen jälkeen alkaa uusi tinkauskierros jossa taas hintaa lasketaan. Kun on päästy loppusuoralle asiat on jo muuttuneet. (Esim. tullut uusia sähkönkuluttajia tai uutta tietoa). "

Add-DesignPrefix "PES sähkön osalta kannattaisi tehdä telakan omalla väellä. 518 tehtiin Kroatiassa jossa oli paljon virheitä ja niitä jouduttiin sitten VAS-suunnittelussa korjaamaan ja paikkaamaan. Pitää selvittää alihankkijan kyky toimittaa taatusti laadukasta työtä jota ei telakka sitten joudu parsimaan kasaan uudestaan." `
                  "Design: PES sähkön osalta kannattaisi tehdä telakan omalla väellä. 518 tehtiin Kroatiassa jossa oli paljon virheitä ja niitä jouduttiin sitten VAS-suunnittelussa korjaamaan ja paikkaamaan. Pitää selvittää alihankkijan kyky toimittaa taatusti laadukasta työtä jota ei telakka sitten joudu parsimaan kasaan uudestaan."

Add-DesignPrefix "Internal communication ok. External communication with suppliers mostly ok." `
                  "Design: Internal communication ok. External communication with suppliers mostly ok."

Add-DesignPrefix "Suunnittelua tehdään samaan aikaan monella osastolla, eli siinä vaiheessa kun telakan pitää jo löydä lukkoon sähkön jakelun keskukset ei vielä ole tarpeeksi lähtötietoja mitä ollaan ostamalla muilla osastoilla." `
                  "Design: Suunnittelua tehdään samaan aikaan monella osastolla, eli siinä vaiheessa kun telakan pitää jo löydä lukkoon sähkön jakelun keskukset ei vielä ole tarpeeksi lähtötietoja mitä ollaan ostamalla muilla osastoilla."

Add-DesignPrefix "Some things went to correct direction but regarding TK the opposite way." `
                  "Design: Some things went to correct direction but regarding TK the opposite way."

Add-DesignPrefix "protolaivat olivat vaikeita tehdä, mutta kyllä niistä opittiin. Osaa asioista ei saatu enää muutettua joten mentiin samalla tavalla kuin edellisissäkin" `
                  "Design: protolaivat olivat vaikeita tehdä, mutta kyllä niistä opittiin. Osaa asioista ei saatu enää muutettua joten mentiin samalla tavalla kuin edellisissäkin"

Write-Output "Done"
